$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 3290.886797766764
$ws.Range("D7").Value = 577.2655423823743
